$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 72 <-> Row 73 swap (everything except column A "id")
$ws.Range("B72").Value = 7517624
$ws.Range("F72").Value = "NK Croatia Dakovo"
$ws.Range("G72").Value = "Valpovka"
$ws.Range("H72").Value = 1
$ws.Range("I72").Value = 6
$ws.Range("K72").Value = 1.4
$ws.Range("L72").Value = 4.8
$ws.Range("M72").Value = 5.25
$ws.Range("N72").Value = 1.4
$ws.Range("O72").Value = 4.8
$ws.Range("P72").Value = 5.25
$ws.Range("Q72").Value = -1.25
$ws.Range("R72").Value = 1.85
$ws.Range("T72").Value = 3.25
$ws.Range("U72").Value = 1.85
$ws.Range("V72").Value = 1.95
$ws.Range("Y72").Value = 4.25
$ws.Range("AB72").Value = 0.8500000000000001

$ws.Range("B73").Value = 7520493
$ws.Range("F73").Value = "NK Medulin"
$ws.Range("G73").Value = "NK Naprijed Hreljin"
$ws.Range("H73").Value = 2
$ws.Range("I73").Value = 4
$ws.Range("K73").Value = 2.8
$ws.Range("L73").Value = 4.1
$ws.Range("M73").Value = 1.95
$ws.Range("N73").Value = 2.8
$ws.Range("O73").Value = 4.1
$ws.Range("P73").Value = 1.95
$ws.Range("Q73").Value = 0.5
$ws.Range("R73").Value = 1.75
$ws.Range("T73").Value = 3
$ws.Range("U73").Value = 1.9
$ws.Range("V73").Value = 1.9
$ws.Range("Y73").Value = 0.95
$ws.Range("AB73").Value = 0.8999999999999999

# Row 101 <-> Row 102 swap (everything except column A "id")
$ws.Range("B101").Value = 8001390
$ws.Range("F101").Value = "Sava Strmec"
$ws.Range("G101").Value = "NK Maksimir"
$ws.Range("I101").Value = 1
$ws.Range("J101").Value = "D"
$ws.Range("N101").Value = 3
$ws.Range("O101").Value = 3.2
$ws.Range("P101").Value = 2.15
$ws.Range("Q101").Value = 0.25
$ws.Range("R101").Value = 1.875
$ws.Range("S101").Value = 1.925
$ws.Range("T101").Value = 2.5
$ws.Range("U101").Value = 1.95
$ws.Range("V101").Value = 1.85
$ws.Range("X101").Value = 2.2
$ws.Range("Y101").Value = -1
$ws.Range("Z101").Value = 0.4375
$ws.Range("AA101").Value = -0.5
$ws.Range("AB101").Value = -1
$ws.Range("AC101").Value = 0.8500000000000001

$ws.Range("B102").Value = 8001389
$ws.Range("F102").Value = "NK Bistra"
$ws.Range("G102").Value = "HNK Segesta"
$ws.Range("I102").Value = 2
$ws.Range("J102").Value = "A"
$ws.Range("N102").Value = 3.8
$ws.Range("O102").Value = 3.5
$ws.Range("P102").Value = 1.75
$ws.Range("Q102").Value = 0.75
$ws.Range("R102").Value = 1.75
$ws.Range("S102").Value = 1.95
$ws.Range("T102").Value = 2.75
$ws.Range("U102").Value = 1.925
$ws.Range("V102").Value = 1.775
$ws.Range("X102").Value = -1
$ws.Range("Y102").Value = 0.75
$ws.Range("Z102").Value = -0.5
$ws.Range("AA102").Value = 0.475
$ws.Range("AB102").Value = 0.4625
$ws.Range("AC102").Value = -0.5
